$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before G - shifts old G..J to H..K and adjusts formulas/col widths.
$ws.Columns("G:G").Insert()

# Match new column G's width to its neighbour (7.2), mirroring the original G-column width.
$ws.Columns("G:G").ColumnWidth = 7.2

# New column G header
$ws.Range("G1").Value = "Occup"

# New column G data: 2 for rows 2-27, 0 for rows 28-42
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 7).Value = 2
}
for ($r = 28; $r -le 42; $r++) {
    $ws.Cells.Item($r, 7).Value = 0
}

Write-Output "done"
